$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (AC1:AE1), copying the header style (bold,
# bordered, centered) from the existing header cell A1 so the new columns
# match the look of the rest of row 1.
$ws.Range("A1").Copy($ws.Range("AC1:AE1"))

$ws.Range("AC1").Value2 = "Wins"
$ws.Range("AD1").Value2 = "Losses"
$ws.Range("AE1").Value2 = "Ties"

# Every data row (2 through 37) gets the same team record: 66 wins, 48
# losses, 1 tie.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 29).Value2 = 66  # AC
    $ws.Cells.Item($r, 30).Value2 = 48  # AD
    $ws.Cells.Item($r, 31).Value2 = 1   # AE
}
